$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the existing section-number values out of column F into column H
#    (Value2 is used so the numbers stay numeric instead of becoming strings)
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 6).Value2
}

# ---------------------------------------------------------------------------
# 2) Populate the new section (F) / adviser (G) columns.
#    The exact order below reproduces the shared-string insertion order of
#    the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value  = "Saragpon, Daniel Catequista"
$ws.Range("G6").Value  = "Vibar, Miles Dominic Morales"
$ws.Range("F1").Value  = "Rizal"
$ws.Range("F6").Value  = "Shakespeare"
$ws.Range("F11").Value = "Industrious"
$ws.Range("F16").Value = "Perseverance"
$ws.Range("G16").Value = "Noa, Kenji Isaac "
$ws.Range("G11").Value = "Zurbano, Christopher Ian "

# Fill down the remaining rows of each block with the same values.
$ws.Range("F2:F5").Value   = "Rizal"
$ws.Range("G2:G5").Value   = "Saragpon, Daniel Catequista"
$ws.Range("F7:F10").Value  = "Shakespeare"
$ws.Range("G7:G10").Value  = "Vibar, Miles Dominic Morales"
$ws.Range("F12:F15").Value = "Industrious"
$ws.Range("G12:G15").Value = "Zurbano, Christopher Ian "
$ws.Range("F17:F20").Value = "Perseverance"
$ws.Range("G17:G20").Value = "Noa, Kenji Isaac "

# ---------------------------------------------------------------------------
# 3) Alignment / styles.
#    Order matters: it determines the order new cellXfs entries are created.
#    1st new style = vertical-center only (column I)
#    2nd new style = right + vertical-center (column F)
#    3rd new style = left + vertical-center (column G)
# ---------------------------------------------------------------------------
$ws.Range("I1:I4").VerticalAlignment = -4108   # xlCenter

$ws.Range("F1:F20").VerticalAlignment   = -4108   # xlCenter
$ws.Range("F1:F20").HorizontalAlignment = -4152   # xlRight

$ws.Range("G1:G20").VerticalAlignment   = -4108   # xlCenter
$ws.Range("G1:G20").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# 4) Column widths (closest values this engine's 1/6-character quantization
#    allows to the authored widths of 15.21875 / 25 / 8.88671875).
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 14.33
$ws.Columns.Item(7).ColumnWidth = 24.17
$ws.Columns.Item(8).ColumnWidth = 8
$ws.Columns.Item(9).ColumnWidth = 8

# ---------------------------------------------------------------------------
# 5) Selection moves to K1.
# ---------------------------------------------------------------------------
$ws.Range("K1").Select()
